$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "trainingaudio/19_papipi1.wav"
$ws.Range("B2").Value = "pngimages/19_burger.png"

$ws.Range("A3").Value = "trainingaudio/08_tipako2.wav"

$ws.Range("A4").Value = "trainingaudio/14_pokoto1.wav"
$ws.Range("B4").Value = "pngimages/14_coffee.png"

$ws.Range("A5").Value = "trainingaudio/17_kotako2.wav"
$ws.Range("B5").Value = "pngimages/17_cracker.png"

$ws.Range("A6").Value = "trainingaudio/22_kakoki1.wav"
$ws.Range("B6").Value = "pngimages/22_egg.png"

$ws.Range("A7").Value = "trainingaudio/01_kitipi1.wav"
$ws.Range("B7").Value = "pngimages/01_gift.png"

$wb.Save()
